$wb = $excel.ActiveWorkbook

# Offense sheet - Home row (row 2) target depth stats update
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 225
$wsOff.Range("C2").Value = 173
$wsOff.Range("D2").Value = 58
$wsOff.Range("E2").Value = 34

# Defense sheet - Home row (row 2) target depth stats update
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 247
$wsDef.Range("C2").Value = 175
$wsDef.Range("D2").Value = 54
$wsDef.Range("E2").Value = 25
$wsDef.Range("F2").Value = 5
